# Update capital structure database:
#  - Row 2 (aggregate row) gets refreshed metrics, company count goes 3 -> 1
#  - Row 3 is replaced with the "EURO Ressources S.A." company data (was Batla Minerals)
#  - Rows 4 and 5 (Euro Ressources duplicate + Auplata Mining Group) are removed entirely

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 2 ----
$ws.Range("B2").Value = 1
$ws.Range("D2").Value = -0.0189
$ws.Range("E2").Value = 0.0319
$ws.Range("G2").Value = 1.015325670498084
$ws.Range("H2").Value = 1.015325670498084
$ws.Range("I2").Value = 0.9578544061302682
$ws.Range("J2").Value = 0.6891230310770541
$ws.Range("K2").Value = 18.2
$ws.Range("L2").Value = 0.6973180076628351
$ws.Range("M2").Value = 14.7
$ws.Range("N2").Value = 0.06223539373412362
$ws.Range("O2").Value = 0.8076923076923077
$ws.Range("P2").Value = 14.7
$ws.Range("Q2").Value = 0.06223539373412362
$ws.Range("R2").Value = 0.8076923076923077
$ws.Range("U2").Value = 35.5
$ws.Range("V2").Value = 0.1502963590177815
$ws.Range("W2").Value = 0.3799582463465553
$ws.Range("X2").Value = 0.05437868308357584
$ws.Range("Y2").Value = 0.3255795632629795
$ws.Range("Z2").Value = 1.601226993865031
$ws.Range("AA2").Value = 1.10344239945467
$ws.Range("AB2").Value = 0.05437868308357584
$ws.Range("AC2").Value = 1.049063716371094
$ws.Range("AD2").Value = 0
$ws.Range("AF2").Value = 0
$ws.Range("AG2").Value = -35.5
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = -0.1768809167912307
$ws.Range("AK2").Value = -2.063953488372093
$ws.Range("AL2").Value = 0
$ws.Range("AM2").Value = -0.472
$ws.Range("AN2").Value = 0
$ws.Range("AO2").ClearContents()
$ws.Range("AP2").Value = -1.397637795275591
$ws.Range("AQ2").Value = -52.96610169491526

# ---- Row 3 (now holds the EURO Ressources S.A. data) ----
$ws.Range("B3").Value = "EURO Ressources S.A. (ENXTPA:EUR)"
$ws.Range("D3").Value = -0.0189
$ws.Range("E3").Value = 0.0319
$ws.Range("G3").Value = 1.015325670498084
$ws.Range("H3").Value = 1.015325670498084
$ws.Range("I3").Value = 0.9578544061302682
$ws.Range("J3").Value = 0.6891230310770541
$ws.Range("K3").Value = 18.2
$ws.Range("L3").Value = 0.6973180076628351
$ws.Range("M3").Value = 14.7
$ws.Range("N3").Value = 0.06223539373412362
$ws.Range("O3").Value = 0.8076923076923077
$ws.Range("P3").Value = 14.7
$ws.Range("Q3").Value = 0.06223539373412362
$ws.Range("R3").Value = 0.8076923076923077
$ws.Range("T3").Value = 0
$ws.Range("U3").Value = 35.5
$ws.Range("V3").Value = 0.1502963590177815
$ws.Range("W3").Value = 0.3799582463465553
$ws.Range("X3").Value = 0.05437868308357584
$ws.Range("Y3").Value = 0.3255795632629795
$ws.Range("Z3").Value = 1.601226993865031
$ws.Range("AA3").Value = 1.10344239945467
$ws.Range("AB3").Value = 0.05437868308357584
$ws.Range("AC3").Value = 1.049063716371094
$ws.Range("AG3").Value = -35.5
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = -0.1768809167912307
$ws.Range("AK3").Value = -2.063953488372093
$ws.Range("AM3").Value = -0.472
$ws.Range("AP3").Value = -1.397637795275591
$ws.Range("AQ3").Value = -52.96610169491526

# ---- Remove rows 4 and 5 (Euro Ressources duplicate row + Auplata Mining Group) ----
$ws.Rows("4:5").Delete()
